$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Bitcoin"
$ws.Range("C2").Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Range("D2").Value = "27.741.66"
$ws.Range("E2").Value = "  +0.40%  "

$ws.Range("B3").Value = "Ethereum"
$ws.Range("C3").Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Range("D3").Value = "1.849.85"
$ws.Range("E3").Value = "  +0.45%  "

$ws.Range("B4").Value = "TetherUSD"
$ws.Range("C4").Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "'314.13"
$ws.Range("E5").Value = "  -0.09%  "

$ws.Range("B6").Value = "USDC"
$ws.Range("C6").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.04%  "

$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").Value = "'0.4320"
$ws.Range("E7").Value = "  +1.38%  "

$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "'0.3660"
$ws.Range("E8").Value = "  +0.41%  "

$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Value = "'45.11"
$ws.Range("E9").Value = "  -1.32%  "

$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.07358"
$ws.Range("E10").Value = "  +1.01%  "

$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").Value = "'0.8781"
$ws.Range("E11").Value = "  -1.97%  "

$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value = "'20.83"
$ws.Range("E12").Value = "  +1.16%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.893.66"
$ws.Range("E13").Value = "  +4.43%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'5.362"
$ws.Range("E14").Value = "  -0.33%  "

$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'6.536"
$ws.Range("E15").Value = "  -0.29%  "

$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").Value = "'0.06972"
$ws.Range("E16").Value = "  +1.66%  "

$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").Value = "'1.000"
$ws.Range("E17").Value = "  -0.28%  "

$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").Value = "'79.93"
$ws.Range("E18").Value = "  +2.03%  "

$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.000009024"
$ws.Range("E19").Value = "  +2.21%  "

$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "'1.002"
$ws.Range("E20").Value = "  +0.25%  "

$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "'15.42"
$ws.Range("E21").Value = "  -0.93%  "

$ws.Range("B22").Value = "WrappedBTC"
$ws.Range("C22").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D22").Value = "27.818.44"
$ws.Range("E22").Value = "  +0.64%  "

$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "'4.982"
$ws.Range("E23").Value = "  +0.15%  "

$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "'10.35"
$ws.Range("E24").Value = "  -2.18%  "

$ws.Range("B25").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C25").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D25").Value = "2.135.10"
$ws.Range("E25").Value = "  +2.74%  "

$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "'1.980"
$ws.Range("E26").Value = "  -3.10%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'156.18"
$ws.Range("E27").Value = "  +1.06%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'18.65"
$ws.Range("E28").Value = "  +2.13%  "

$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").Value = "'120.16"
$ws.Range("E29").Value = "  +7.82%  "

$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "'5.262"
$ws.Range("E30").Value = "  +0.55%  "

$ws.Range("B31").Value = "LidoDAOToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D31").Value = "'1.886"
$ws.Range("E31").Value = "  +2.95%  "

$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "'0.08911"
$ws.Range("E32").Value = "  +0.25%  "

$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'0.7616"
$ws.Range("E33").Value = "  -1.68%  "

$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'4.566"
$ws.Range("E34").Value = "  -0.03%  "

$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'2.964"
$ws.Range("E35").Value = "  -1.04%  "

$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'1.134"
$ws.Range("E36").Value = "  +3.28%  "

$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").Value = "'1.000"
$ws.Range("E37").Value = "  +0.11%  "

$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'1.111"
$ws.Range("E38").Value = "  +1.20%  "

$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.05446"
$ws.Range("E39").Value = "  +0.15%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.01938"
$ws.Range("E40").Value = "  +0.59%  "

$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "'2.826"
$ws.Range("E41").Value = "  +1.98%  "

$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "'0.5099"
$ws.Range("E42").Value = "  +0.74%  "

$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "'0.1663"
$ws.Range("E43").Value = "  +1.17%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'6.613"
$ws.Range("E44").Value = "  -2.77%  "

$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "'8.419"
$ws.Range("E45").Value = "  +2.41%  "

$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").Value = "'0.06547"
$ws.Range("E46").Value = "  -1.36%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'10.40"
$ws.Range("E47").Value = "  -0.19%  "

$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "'105.37"
$ws.Range("E48").Value = "  -0.46%  "

$ws.Range("B49").Value = "Decentraland"
$ws.Range("C49").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D49").Value = "'0.4674"
$ws.Range("E49").Value = "  -0.81%  "

$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").Value = "'0.9998"
$ws.Range("E50").Value = "  -0.05%  "

$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "'1.639"
$ws.Range("E51").Value = "  +0.21%  "
